$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "...among fans of many different sports today..."
#    -> "...among fans of a variety of sports today..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("many different ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "a variety of ", 2)

# ---------------------------------------------------------------------------
# 2. Typo fix: "...National Football League (NFL) though data mining..."
#    -> "...National Football League (NFL) through data mining..."
#    (split the run the way Word does when you insert a single "r" in place,
#     i.e. "...(NFL) th" | "r" | "ough data mining techniques.")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("though data mining", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "through data mining", 2)

$probe = $d.Content
$probe.Find.Execute("(NFL) th", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
if ($probe.Find.Found) {
    $p1 = $probe.End

    $probe2 = $d.Content
    $probe2.Find.Execute("(NFL) thr", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
    $p2 = $probe2.End

    $probe3 = $d.Content
    $probe3.Find.Execute("ough data mining techniques.", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
    $p3 = $probe3.End

    $rA = $d.Range($p1, $p2)
    $rA.Font.Bold = $true
    $rA.Font.Bold = $false

    $rB = $d.Range($p2, $p3)
    $rB.Font.Bold = $true
    $rB.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 3. "-Sep: Initial data accumulation" -> "-Sep: Initial data acquisition"
#    (split so "Initial data a" stays its own run, like Word does when you
#     retype "cquisition" over "ccumulation").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Initial data accumulation", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Initial data acquisition", 2)

$probe4 = $d.Content
$probe4.Find.Execute("Initial data a", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
if ($probe4.Find.Found) {
    $q1 = $probe4.End

    $probe5 = $d.Content
    $probe5.Find.Execute("cquisition", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
    $q2 = $probe5.End

    $rC = $d.Range($q1, $q2)
    $rC.Font.Bold = $true
    $rC.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 4. Move the (hidden, singleton) "_GoBack" bookmark -- Word stamps it at the
#    location of the author's last edit.  It used to sit right after the
#    "Project Objectives:" tab; the last thing Kevin touched was the
#    footballlocks.com reference link, so it now belongs between "co" and
#    "m" of that URL.  Word can't normally split a hyperlink run through the
#    object model, so nudge it open with a throw-away character, drop the
#    bookmark in the gap, then remove the throw-away character again.
# ---------------------------------------------------------------------------
$linkProbe = $d.Content
$linkProbe.Find.Execute("http://www.footballlocks.co", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
if ($linkProbe.Find.Found) {
    $splitPos = $linkProbe.End

    $gap = $d.Range($splitPos, $splitPos)
    $gap.InsertAfter("X")

    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $throwAway = $d.Range($splitPos + 1, $splitPos + 2)
    $throwAway.Delete()
}
